$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D3").Value = -7.984
$ws.Range("C7").Value = -12.844
$ws.Range("A10").Value = -21.831
$ws.Range("E10").Value = 16.466
$ws.Range("A12").Value = -21.727
$ws.Range("E14").Value = 17.027
$ws.Range("C15").Value = -14.056
$ws.Range("A18").Value = -22.265
$ws.Range("D18").Value = -8.84
$ws.Range("D19").Value = -8.113
$ws.Range("C20").Value = -12.673
$ws.Range("D27").Value = -8.624000000000001
$ws.Range("C29").Value = -12.124
$ws.Range("C30").Value = -12.989
$ws.Range("C31").Value = -13.439
$ws.Range("E32").Value = 16.876
$ws.Range("E35").Value = 16.358
$ws.Range("A37").Value = -20.029
$ws.Range("C40").Value = -12.782
$ws.Range("D42").Value = -8.258000000000001
$ws.Range("E43").Value = 17.138
$ws.Range("D44").Value = -7.529999999999999
$ws.Range("D47").Value = -7.566999999999998
$ws.Range("E49").Value = 16.349
$ws.Range("A55").Value = -21.894
$ws.Range("E56").Value = 16.166
$ws.Range("D58").Value = -8.400999999999998
$ws.Range("A68").Value = -21.707
$ws.Range("C68").Value = -11.001
$ws.Range("E69").Value = 17.438
$ws.Range("D73").Value = -7.877000000000001
$ws.Range("C76").Value = -13.045
$ws.Range("A77").Value = -20.843
$ws.Range("A78").Value = -19.951
$ws.Range("E81").Value = 16.558
$ws.Range("C87").Value = -13.458
$ws.Range("C88").Value = -13.142
$ws.Range("E92").Value = 17.939
$ws.Range("D95").Value = -7.567
$ws.Range("C96").Value = -12.644
$ws.Range("C98").Value = -13.645
$ws.Range("C101").Value = -13.048
$ws.Range("D101").Value = -8.061999999999999
$ws.Range("C102").Value = -13.086
